# Atualização automática dos dados: Tue Jan 20 09:38:08 UTC 2026
# Applies the updated "Entrada" dashboard figures to the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entrada")

# Row 2 - DEVOLUÇÃO
$ws.Range("B2").Value = "R$ 419.138,96"
$ws.Range("D2").Value = "R$ 419.138,96"
$ws.Range("E2").Value = "R$ 419.138,96"

# Row 3 - FERRAMENTAS/ MATRIZARIA
$ws.Range("B3").Value = "R$ 399.858,72"
$ws.Range("D3").Value = "R$ 399.858,72"
$ws.Range("F3").Value = "54,78 %"

# Row 4 - MATERIA PRIMA
$ws.Range("B4").Value = "R$ 289.559,50"
$ws.Range("D4").Value = "R$ 289.559,50"
$ws.Range("F4").Value = "28.955.950,00 %"

# Row 5 - REFUGO REAL (PROCESSO)
$ws.Range("B5").Value = "R$ 271.490,35"
$ws.Range("D5").Value = "R$ 271.490,35"
$ws.Range("E5").Value = "R$ 271.490,35"

# Row 7 - label swapped from REFUGO MP+CP* to MANUTENCAO
$ws.Range("A7").Value = "MANUTENCAO"
$ws.Range("B7").Value = "R$ 179.839,62"
$ws.Range("C7").Value = "R$ 276.953,04"
$ws.Range("D7").Value = "R$ 456.792,66"
$ws.Range("E7").Value = "R$ 480.000,00"
$ws.Range("F7").Value = "95,17 %"

# Row 8 - label swapped from MANUTENCAO to REFUGO MP+CP*
$ws.Range("A8").Value = "REFUGO MP+CP*"
$ws.Range("B8").Value = "R$ 159.821,42"
$ws.Range("C8").Value = "R$ 0,00"
$ws.Range("D8").Value = "R$ 159.821,42"
$ws.Range("E8").Value = "R$ 285.000,00"
$ws.Range("F8").Value = "56,08 %"

# Row 9 - label swapped from OLEOS E LUBRIFICANTES to CUSTO DESENVOLVIMENTO
$ws.Range("A9").Value = "CUSTO DESENVOLVIMENTO"
$ws.Range("B9").Value = "R$ 148.010,88"
$ws.Range("C9").Value = "R$ 0,00"
$ws.Range("D9").Value = "R$ 148.010,88"
$ws.Range("E9").Value = "R$ 148.010,88"
$ws.Range("F9").Value = "100,00 %"

# Row 10 - label swapped from CUSTO DESENVOLVIMENTO to OLEOS E LUBRIFICANTES
$ws.Range("A10").Value = "OLEOS E LUBRIFICANTES"
$ws.Range("B10").Value = "R$ 82.091,98"
$ws.Range("C10").Value = "R$ 109.354,76"
$ws.Range("D10").Value = "R$ 191.446,74"
$ws.Range("E10").Value = "R$ 280.000,00"
$ws.Range("F10").Value = "68,37 %"

# Row 11 - DESP. INDUSTRIAL
$ws.Range("B11").Value = "R$ 57.567,20"
$ws.Range("C11").Value = "R$ 103.938,96"
$ws.Range("D11").Value = "R$ 161.506,16"
$ws.Range("F11").Value = "34,36 %"

# Row 12 - EMBALAGENS
$ws.Range("C12").Value = "R$ 64.920,59"
$ws.Range("D12").Value = "R$ 112.894,67"
$ws.Range("F12").Value = "70,56 %"

# Row 14 - label swapped from CUSTO COM DESENVOLVIMENTO to FERRAMENTARIA/MAN FR
$ws.Range("A14").Value = "FERRAMENTARIA/MAN FR"
$ws.Range("B14").Value = "R$ 8.619,54"
$ws.Range("D14").Value = "R$ 8.619,54"
$ws.Range("E14").Value = "R$ 35.000,00"
$ws.Range("F14").Value = "24,63 %"

# Row 15 - label swapped from FERRAMENTARIA/MAN FR to CUSTO COM DESENVOLVIMENTO
$ws.Range("A15").Value = "CUSTO COM DESENVOLVIMENTO"
$ws.Range("B15").Value = "R$ 8.301,08"
$ws.Range("D15").Value = "R$ 8.301,08"
$ws.Range("E15").Value = "R$ 8.301,08"
$ws.Range("F15").Value = "100,00 %"

# Row 18 - Total Geral
$ws.Range("B18").Value = "R$ 2.141.695,10"
$ws.Range("C18").Value = "R$ 559.201,15"
$ws.Range("D18").Value = "R$ 2.700.896,26"
$ws.Range("E18").Value = "R$ 3.986.942,27"
$ws.Range("F18").Value = "67,74 %"
